$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")
    $val2 = $cellRow2.Value2
    $val3 = $cellRow3.Value2
    $cellRow2.Value2 = $val3
    $cellRow3.Value2 = $val2
}
